$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44231
$ws.Cells.Item(2, 8).Value = "Perfection"
$ws.Cells.Item(2, 10).Value = 110
$ws.Cells.Item(2, 11).Value = 20000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 20000
$ws.Cells.Item(2, 16).Value = 800
# Row 3
$ws.Cells.Item(3, 4).Value = 44539
$ws.Cells.Item(3, 10).Value = 225
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 15000
$ws.Cells.Item(3, 13).Value = 15000
$ws.Cells.Item(3, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(3, 16).Value = 600
# Row 4
$ws.Cells.Item(4, 4).Value = 44585
$ws.Cells.Item(4, 10).Value = 65
$ws.Cells.Item(4, 11).Value = 26000
$ws.Cells.Item(4, 12).Value = 26000
$ws.Cells.Item(4, 13).Value = 26000
$ws.Cells.Item(4, 16).Value = 1040
# Row 5
$ws.Cells.Item(5, 4).Value = 44469
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 34000
$ws.Cells.Item(5, 12).Value = 34000
$ws.Cells.Item(5, 13).Value = 34000
$ws.Cells.Item(5, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 1360
# Row 6
$ws.Cells.Item(6, 4).Value = 44875
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(6, 11).Value = 25000
$ws.Cells.Item(6, 12).Value = 25000
$ws.Cells.Item(6, 13).Value = 25000
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 1000
# Row 7
$ws.Cells.Item(7, 4).Value = 44483
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 25000
$ws.Cells.Item(7, 12).Value = 25000
$ws.Cells.Item(7, 13).Value = 25000
$ws.Cells.Item(7, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1000
# Row 8
$ws.Cells.Item(8, 4).Value = 44162
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 23000
$ws.Cells.Item(8, 12).Value = 23000
$ws.Cells.Item(8, 13).Value = 23000
$ws.Cells.Item(8, 16).Value = 920
# Row 9
$ws.Cells.Item(9, 4).Value = 44162
$ws.Cells.Item(9, 10).Value = 40
$ws.Cells.Item(9, 11).Value = 21000
$ws.Cells.Item(9, 12).Value = 21000
$ws.Cells.Item(9, 13).Value = 21000
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 840
# Row 10
$ws.Cells.Item(10, 4).Value = 44165
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 22000
$ws.Cells.Item(10, 12).Value = 22000
$ws.Cells.Item(10, 13).Value = 22000
$ws.Cells.Item(10, 16).Value = 880
# Row 11
$ws.Cells.Item(11, 4).Value = 44210
$ws.Cells.Item(11, 10).Value = 150
$ws.Cells.Item(11, 11).Value = 17000
$ws.Cells.Item(11, 12).Value = 17000
$ws.Cells.Item(11, 13).Value = 17000
$ws.Cells.Item(11, 16).Value = 680
# Row 12
$ws.Cells.Item(12, 4).Value = 44554
$ws.Cells.Item(12, 10).Value = 75
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14667
$ws.Cells.Item(12, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(12, 16).Value = 587
# Row 13
$ws.Cells.Item(13, 4).Value = 44959
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 35000
$ws.Cells.Item(13, 12).Value = 35000
$ws.Cells.Item(13, 13).Value = 35000
$ws.Cells.Item(13, 16).Value = 1400
# Row 14
$ws.Cells.Item(14, 4).Value = 44935
$ws.Cells.Item(14, 10).Value = 95
$ws.Cells.Item(14, 11).Value = 32000
$ws.Cells.Item(14, 12).Value = 33000
$ws.Cells.Item(14, 13).Value = 32421
$ws.Cells.Item(14, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(14, 16).Value = 1297
# Row 15
$ws.Cells.Item(15, 4).Value = 44571
$ws.Cells.Item(15, 10).Value = 95
$ws.Cells.Item(15, 11).Value = 18000
$ws.Cells.Item(15, 12).Value = 18000
$ws.Cells.Item(15, 13).Value = 18000
$ws.Cells.Item(15, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(15, 16).Value = 720
# Row 16
$ws.Cells.Item(16, 4).Value = 44921
$ws.Cells.Item(16, 10).Value = 55
$ws.Cells.Item(16, 11).Value = 25000
$ws.Cells.Item(16, 12).Value = 25000
$ws.Cells.Item(16, 13).Value = 25000
$ws.Cells.Item(16, 16).Value = 1000
# Row 17
$ws.Cells.Item(17, 4).Value = 44224
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 10).Value = 320
$ws.Cells.Item(17, 11).Value = 17000
$ws.Cells.Item(17, 12).Value = 17000
$ws.Cells.Item(17, 13).Value = 17000
$ws.Cells.Item(17, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(17, 16).Value = 680
# Row 18
$ws.Cells.Item(18, 4).Value = 44497
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 17000
$ws.Cells.Item(18, 12).Value = 20000
$ws.Cells.Item(18, 13).Value = 18500
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 740
# Row 19
$ws.Cells.Item(19, 4).Value = 44533
$ws.Cells.Item(19, 10).Value = 110
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 15000
$ws.Cells.Item(19, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(19, 16).Value = 600
# Row 20
$ws.Cells.Item(20, 4).Value = 44533
$ws.Cells.Item(20, 10).Value = 95
$ws.Cells.Item(20, 15).Value = "Región del Maule"
# Row 21
$ws.Cells.Item(21, 4).Value = 44232
$ws.Cells.Item(21, 10).Value = 110
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 13).Value = 20000
$ws.Cells.Item(21, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(21, 16).Value = 800
# Row 22
$ws.Cells.Item(22, 4).Value = 44603
$ws.Cells.Item(22, 10).Value = 185
$ws.Cells.Item(22, 11).Value = 22000
$ws.Cells.Item(22, 12).Value = 22000
$ws.Cells.Item(22, 13).Value = 22000
$ws.Cells.Item(22, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(22, 16).Value = 880
# Row 23
$ws.Cells.Item(23, 4).Value = 44487
$ws.Cells.Item(23, 10).Value = 110
$ws.Cells.Item(23, 11).Value = 25000
$ws.Cells.Item(23, 12).Value = 25000
$ws.Cells.Item(23, 13).Value = 25000
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 1000
# Row 24
$ws.Cells.Item(24, 4).Value = 45006
$ws.Cells.Item(24, 10).Value = 25
$ws.Cells.Item(24, 11).Value = 33000
$ws.Cells.Item(24, 12).Value = 33000
$ws.Cells.Item(24, 13).Value = 33000
$ws.Cells.Item(24, 16).Value = 1320
# Row 25
$ws.Cells.Item(25, 4).Value = 44588
$ws.Cells.Item(25, 10).Value = 65
$ws.Cells.Item(25, 11).Value = 24000
$ws.Cells.Item(25, 12).Value = 24000
$ws.Cells.Item(25, 13).Value = 24000
$ws.Cells.Item(25, 16).Value = 960
# Row 26
$ws.Cells.Item(26, 4).Value = 44581
$ws.Cells.Item(26, 10).Value = 40
$ws.Cells.Item(26, 11).Value = 26000
$ws.Cells.Item(26, 12).Value = 26000
$ws.Cells.Item(26, 13).Value = 26000
$ws.Cells.Item(26, 16).Value = 1040
# Row 27
$ws.Cells.Item(27, 4).Value = 44172
$ws.Cells.Item(27, 10).Value = 85
$ws.Cells.Item(27, 11).Value = 19000
$ws.Cells.Item(27, 12).Value = 20000
$ws.Cells.Item(27, 13).Value = 19412
$ws.Cells.Item(27, 16).Value = 776
# Row 28
$ws.Cells.Item(28, 4).Value = 44970
$ws.Cells.Item(28, 10).Value = 70
$ws.Cells.Item(28, 11).Value = 25000
$ws.Cells.Item(28, 12).Value = 30000
$ws.Cells.Item(28, 13).Value = 27143
$ws.Cells.Item(28, 16).Value = 1086
# Row 29
$ws.Cells.Item(29, 4).Value = 44203
$ws.Cells.Item(29, 10).Value = 160
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 16875
$ws.Cells.Item(29, 16).Value = 675
# Row 30
$ws.Cells.Item(30, 4).Value = 44208
$ws.Cells.Item(30, 10).Value = 55
$ws.Cells.Item(30, 11).Value = 17000
$ws.Cells.Item(30, 12).Value = 17000
$ws.Cells.Item(30, 13).Value = 17000
$ws.Cells.Item(30, 16).Value = 680
# Row 31
$ws.Cells.Item(31, 4).Value = 44551
$ws.Cells.Item(31, 10).Value = 155
$ws.Cells.Item(31, 11).Value = 15000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = 15000
$ws.Cells.Item(31, 16).Value = 600
# Row 32
$ws.Cells.Item(32, 4).Value = 44565
$ws.Cells.Item(32, 10).Value = 20
$ws.Cells.Item(32, 11).Value = 26000
$ws.Cells.Item(32, 12).Value = 26000
$ws.Cells.Item(32, 13).Value = 26000
$ws.Cells.Item(32, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(32, 16).Value = 1040
# Row 33
$ws.Cells.Item(33, 4).Value = 44505
$ws.Cells.Item(33, 10).Value = 125
$ws.Cells.Item(33, 11).Value = 16000
$ws.Cells.Item(33, 12).Value = 16000
$ws.Cells.Item(33, 13).Value = 16000
$ws.Cells.Item(33, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(33, 16).Value = 640
# Row 34
$ws.Cells.Item(34, 4).Value = 44159
$ws.Cells.Item(34, 10).Value = 100
$ws.Cells.Item(34, 11).Value = 20000
$ws.Cells.Item(34, 12).Value = 20000
$ws.Cells.Item(34, 13).Value = 20000
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 800
# Row 35
$ws.Cells.Item(35, 4).Value = 44925
$ws.Cells.Item(35, 10).Value = 120
$ws.Cells.Item(35, 11).Value = 25000
$ws.Cells.Item(35, 12).Value = 25000
$ws.Cells.Item(35, 13).Value = 25000
$ws.Cells.Item(35, 16).Value = 1000
# Row 36
$ws.Cells.Item(36, 4).Value = 44168
$ws.Cells.Item(36, 10).Value = 75
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 20000
$ws.Cells.Item(36, 16).Value = 800
# Row 37
$ws.Cells.Item(37, 4).Value = 44195
$ws.Cells.Item(37, 10).Value = 110
$ws.Cells.Item(37, 11).Value = 25000
$ws.Cells.Item(37, 12).Value = 25000
$ws.Cells.Item(37, 13).Value = 25000
$ws.Cells.Item(37, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(37, 16).Value = 1000
# Row 38
$ws.Cells.Item(38, 4).Value = 44204
$ws.Cells.Item(38, 10).Value = 40
$ws.Cells.Item(38, 11).Value = 15000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = 15000
$ws.Cells.Item(38, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(38, 16).Value = 600
# Row 39
$ws.Cells.Item(39, 4).Value = 44503
$ws.Cells.Item(39, 10).Value = 75
$ws.Cells.Item(39, 11).Value = 15000
$ws.Cells.Item(39, 12).Value = 15000
$ws.Cells.Item(39, 13).Value = 15000
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 600
# Row 40
$ws.Cells.Item(40, 4).Value = 44858
$ws.Cells.Item(40, 10).Value = 215
$ws.Cells.Item(40, 11).Value = 18000
$ws.Cells.Item(40, 12).Value = 20000
$ws.Cells.Item(40, 13).Value = 19163
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 767
# Row 41
$ws.Cells.Item(41, 4).Value = 44942
$ws.Cells.Item(41, 10).Value = 45
$ws.Cells.Item(41, 11).Value = 30000
$ws.Cells.Item(41, 12).Value = 32000
$ws.Cells.Item(41, 13).Value = 30667
$ws.Cells.Item(41, 16).Value = 1227
# Row 42
$ws.Cells.Item(42, 4).Value = 44536
$ws.Cells.Item(42, 10).Value = 290
$ws.Cells.Item(42, 11).Value = 13000
$ws.Cells.Item(42, 13).Value = 14138
$ws.Cells.Item(42, 16).Value = 566
# Row 43
$ws.Cells.Item(43, 4).Value = 44559
$ws.Cells.Item(43, 10).Value = 110
# Row 44
$ws.Cells.Item(44, 4).Value = 44964
$ws.Cells.Item(44, 10).Value = 80
$ws.Cells.Item(44, 11).Value = 25000
$ws.Cells.Item(44, 12).Value = 25000
$ws.Cells.Item(44, 13).Value = 25000
$ws.Cells.Item(44, 16).Value = 1000
# Row 45
$ws.Cells.Item(45, 4).Value = 44923
$ws.Cells.Item(45, 10).Value = 115
$ws.Cells.Item(45, 11).Value = 25000
$ws.Cells.Item(45, 12).Value = 25000
$ws.Cells.Item(45, 13).Value = 25000
$ws.Cells.Item(45, 16).Value = 1000
# Row 46
$ws.Cells.Item(46, 4).Value = 44923
$ws.Cells.Item(46, 10).Value = 115
$ws.Cells.Item(46, 11).Value = 25000
$ws.Cells.Item(46, 12).Value = 25000
$ws.Cells.Item(46, 13).Value = 25000
$ws.Cells.Item(46, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(46, 16).Value = 1000
# Row 47
$ws.Cells.Item(47, 4).Value = 44945
$ws.Cells.Item(47, 8).Value = "Perfection"
$ws.Cells.Item(47, 10).Value = 80
$ws.Cells.Item(47, 11).Value = 30000
$ws.Cells.Item(47, 12).Value = 30000
$ws.Cells.Item(47, 13).Value = 30000
$ws.Cells.Item(47, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(47, 16).Value = 1200
# Row 48
$ws.Cells.Item(48, 4).Value = 44239
$ws.Cells.Item(48, 10).Value = 110
$ws.Cells.Item(48, 11).Value = 17000
$ws.Cells.Item(48, 12).Value = 17000
$ws.Cells.Item(48, 13).Value = 17000
$ws.Cells.Item(48, 16).Value = 680
# Row 49
$ws.Cells.Item(49, 4).Value = 44243
$ws.Cells.Item(49, 10).Value = 95
$ws.Cells.Item(49, 11).Value = 17000
$ws.Cells.Item(49, 12).Value = 17000
$ws.Cells.Item(49, 13).Value = 17000
$ws.Cells.Item(49, 16).Value = 680
# Row 50
$ws.Cells.Item(50, 4).Value = 44271
$ws.Cells.Item(50, 10).Value = 80
# Row 51
$ws.Cells.Item(51, 4).Value = 44613
$ws.Cells.Item(51, 10).Value = 30
$ws.Cells.Item(51, 11).Value = 26000
$ws.Cells.Item(51, 12).Value = 26000
$ws.Cells.Item(51, 13).Value = 26000
$ws.Cells.Item(51, 16).Value = 1040
# Row 52
$ws.Cells.Item(52, 4).Value = 44508
$ws.Cells.Item(52, 10).Value = 30
$ws.Cells.Item(52, 15).Value = "Región Metropolitana"
# Row 53
$ws.Cells.Item(53, 4).Value = 44508
$ws.Cells.Item(53, 10).Value = 70
$ws.Cells.Item(53, 12).Value = 17000
$ws.Cells.Item(53, 13).Value = 16571
$ws.Cells.Item(53, 15).Value = "Región del Maule"
$ws.Cells.Item(53, 16).Value = 663
# Row 54
$ws.Cells.Item(54, 4).Value = 44494
$ws.Cells.Item(54, 10).Value = 30
$ws.Cells.Item(54, 11).Value = 21000
$ws.Cells.Item(54, 12).Value = 21000
$ws.Cells.Item(54, 13).Value = 21000
$ws.Cells.Item(54, 15).Value = "Región Metropolitana"
$ws.Cells.Item(54, 16).Value = 840
# Row 55
$ws.Cells.Item(55, 4).Value = 45015
$ws.Cells.Item(55, 10).Value = 40
$ws.Cells.Item(55, 11).Value = 23000
$ws.Cells.Item(55, 12).Value = 23000
$ws.Cells.Item(55, 13).Value = 23000
$ws.Cells.Item(55, 16).Value = 920
# Row 56
$ws.Cells.Item(56, 4).Value = 44179
$ws.Cells.Item(56, 11).Value = 17000
$ws.Cells.Item(56, 12).Value = 17000
$ws.Cells.Item(56, 13).Value = 17000
$ws.Cells.Item(56, 16).Value = 680
# Row 57
$ws.Cells.Item(57, 4).Value = 44879
$ws.Cells.Item(57, 10).Value = 160
$ws.Cells.Item(57, 11).Value = 22000
$ws.Cells.Item(57, 12).Value = 23000
$ws.Cells.Item(57, 13).Value = 22500
$ws.Cells.Item(57, 15).Value = "Región Metropolitana"
$ws.Cells.Item(57, 16).Value = 900
# Row 58
$ws.Cells.Item(58, 4).Value = 44915
$ws.Cells.Item(58, 10).Value = 80
$ws.Cells.Item(58, 11).Value = 25000
$ws.Cells.Item(58, 12).Value = 25000
$ws.Cells.Item(58, 13).Value = 25000
$ws.Cells.Item(58, 16).Value = 1000
# Row 59
$ws.Cells.Item(59, 4).Value = 44880
$ws.Cells.Item(59, 10).Value = 40
$ws.Cells.Item(59, 11).Value = 20000
$ws.Cells.Item(59, 12).Value = 20000
$ws.Cells.Item(59, 13).Value = 20000
$ws.Cells.Item(59, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(59, 15).Value = "Región Metropolitana"
$ws.Cells.Item(59, 16).Value = 800
# Row 60
$ws.Cells.Item(60, 4).Value = 44515
$ws.Cells.Item(60, 10).Value = 115
$ws.Cells.Item(60, 11).Value = 16000
$ws.Cells.Item(60, 12).Value = 16000
$ws.Cells.Item(60, 13).Value = 16000
$ws.Cells.Item(60, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(60, 15).Value = "Región del Maule"
$ws.Cells.Item(60, 16).Value = 640
# Row 61
$ws.Cells.Item(61, 4).Value = 44546
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 10).Value = 95
$ws.Cells.Item(61, 11).Value = 15000
$ws.Cells.Item(61, 12).Value = 15000
$ws.Cells.Item(61, 13).Value = 15000
$ws.Cells.Item(61, 16).Value = 600
# Row 62
$ws.Cells.Item(62, 4).Value = 44936
$ws.Cells.Item(62, 10).Value = 55
$ws.Cells.Item(62, 11).Value = 32000
$ws.Cells.Item(62, 12).Value = 32000
$ws.Cells.Item(62, 13).Value = 32000
$ws.Cells.Item(62, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(62, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(62, 16).Value = 1280
# Row 63
$ws.Cells.Item(63, 4).Value = 44550
$ws.Cells.Item(63, 10).Value = 95
$ws.Cells.Item(63, 11).Value = 15000
$ws.Cells.Item(63, 12).Value = 15000
$ws.Cells.Item(63, 13).Value = 15000
$ws.Cells.Item(63, 16).Value = 600
# Row 64
$ws.Cells.Item(64, 4).Value = 44225
$ws.Cells.Item(64, 10).Value = 25
$ws.Cells.Item(64, 11).Value = 17000
$ws.Cells.Item(64, 12).Value = 17000
$ws.Cells.Item(64, 13).Value = 17000
$ws.Cells.Item(64, 16).Value = 680
# Row 65
$ws.Cells.Item(65, 4).Value = 44176
$ws.Cells.Item(65, 10).Value = 20
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 20000
$ws.Cells.Item(65, 13).Value = 20000
$ws.Cells.Item(65, 16).Value = 800
# Row 66
$ws.Cells.Item(66, 4).Value = 44868
$ws.Cells.Item(66, 10).Value = 170
$ws.Cells.Item(66, 11).Value = 20000
$ws.Cells.Item(66, 12).Value = 22000
$ws.Cells.Item(66, 13).Value = 21059
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value = 842
# Row 67
$ws.Cells.Item(67, 4).Value = 44973
$ws.Cells.Item(67, 10).Value = 40
$ws.Cells.Item(67, 11).Value = 30000
$ws.Cells.Item(67, 12).Value = 30000
$ws.Cells.Item(67, 13).Value = 30000
$ws.Cells.Item(67, 16).Value = 1200
# Row 68
$ws.Cells.Item(68, 4).Value = 44188
$ws.Cells.Item(68, 10).Value = 50
$ws.Cells.Item(68, 11).Value = 25000
$ws.Cells.Item(68, 12).Value = 25000
$ws.Cells.Item(68, 13).Value = 25000
$ws.Cells.Item(68, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(68, 16).Value = 1000
# Row 69
$ws.Cells.Item(69, 4).Value = 44236
$ws.Cells.Item(69, 10).Value = 95
# Row 70
$ws.Cells.Item(70, 4).Value = 44221
$ws.Cells.Item(70, 10).Value = 210
$ws.Cells.Item(70, 11).Value = 18000
$ws.Cells.Item(70, 12).Value = 18000
$ws.Cells.Item(70, 13).Value = 18000
$ws.Cells.Item(70, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(70, 16).Value = 720
# Row 71
$ws.Cells.Item(71, 4).Value = 44511
$ws.Cells.Item(71, 10).Value = 700
$ws.Cells.Item(71, 11).Value = 16000
$ws.Cells.Item(71, 12).Value = 17000
$ws.Cells.Item(71, 13).Value = 16571
$ws.Cells.Item(71, 15).Value = "Región del Maule"
$ws.Cells.Item(71, 16).Value = 663
# Row 72
$ws.Cells.Item(72, 4).Value = 44175
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 10).Value = 40
$ws.Cells.Item(72, 11).Value = 20000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 20000
$ws.Cells.Item(72, 16).Value = 800
# Row 73
$ws.Cells.Item(73, 4).Value = 44512
$ws.Cells.Item(73, 8).Value = "Perfection"
$ws.Cells.Item(73, 10).Value = 50
$ws.Cells.Item(73, 15).Value = "Región del Maule"
# Row 74
$ws.Cells.Item(74, 4).Value = 44235
$ws.Cells.Item(74, 10).Value = 250
$ws.Cells.Item(74, 11).Value = 17000
$ws.Cells.Item(74, 12).Value = 17000
$ws.Cells.Item(74, 13).Value = 17000
$ws.Cells.Item(74, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(74, 16).Value = 680
# Row 75
$ws.Cells.Item(75, 4).Value = 44526
$ws.Cells.Item(75, 10).Value = 40
$ws.Cells.Item(75, 11).Value = 15000
$ws.Cells.Item(75, 12).Value = 15000
$ws.Cells.Item(75, 13).Value = 15000
$ws.Cells.Item(75, 16).Value = 600
# Row 76
$ws.Cells.Item(76, 4).Value = 44526
$ws.Cells.Item(76, 10).Value = 40
$ws.Cells.Item(76, 11).Value = 15000
$ws.Cells.Item(76, 12).Value = 15000
$ws.Cells.Item(76, 13).Value = 15000
$ws.Cells.Item(76, 15).Value = "Región del Maule"
$ws.Cells.Item(76, 16).Value = 600
# Row 77
$ws.Cells.Item(77, 4).Value = 44169
$ws.Cells.Item(77, 10).Value = 135
$ws.Cells.Item(77, 11).Value = 20000
$ws.Cells.Item(77, 12).Value = 22000
$ws.Cells.Item(77, 13).Value = 20593
$ws.Cells.Item(77, 16).Value = 824
# Row 78
$ws.Cells.Item(78, 4).Value = 44160
$ws.Cells.Item(78, 10).Value = 80
$ws.Cells.Item(78, 11).Value = 21000
$ws.Cells.Item(78, 12).Value = 22000
$ws.Cells.Item(78, 13).Value = 21625
$ws.Cells.Item(78, 15).Value = "Región del Maule"
$ws.Cells.Item(78, 16).Value = 865
# Row 79
$ws.Cells.Item(79, 4).Value = 44166
$ws.Cells.Item(79, 11).Value = 21000
$ws.Cells.Item(79, 12).Value = 22000
$ws.Cells.Item(79, 13).Value = 21450
$ws.Cells.Item(79, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(79, 16).Value = 858
# Row 80
$ws.Cells.Item(80, 4).Value = 44273
$ws.Cells.Item(80, 11).Value = 20000
$ws.Cells.Item(80, 12).Value = 20000
$ws.Cells.Item(80, 13).Value = 20000
$ws.Cells.Item(80, 16).Value = 800
# Row 81
$ws.Cells.Item(81, 4).Value = 44882
$ws.Cells.Item(81, 10).Value = 80
$ws.Cells.Item(81, 11).Value = 16000
$ws.Cells.Item(81, 12).Value = 16000
$ws.Cells.Item(81, 13).Value = 16000
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 640
# Row 82
$ws.Cells.Item(82, 4).Value = 44186
$ws.Cells.Item(82, 10).Value = 40
$ws.Cells.Item(82, 11).Value = 20000
$ws.Cells.Item(82, 12).Value = 20000
$ws.Cells.Item(82, 13).Value = 20000
$ws.Cells.Item(82, 16).Value = 800
# Row 83
$ws.Cells.Item(83, 4).Value = 44917
$ws.Cells.Item(83, 10).Value = 65
$ws.Cells.Item(83, 11).Value = 25000
$ws.Cells.Item(83, 12).Value = 25000
$ws.Cells.Item(83, 13).Value = 25000
$ws.Cells.Item(83, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(83, 16).Value = 1000
# Row 84
$ws.Cells.Item(84, 4).Value = 44566
$ws.Cells.Item(84, 10).Value = 20
$ws.Cells.Item(84, 11).Value = 26000
$ws.Cells.Item(84, 12).Value = 26000
$ws.Cells.Item(84, 13).Value = 26000
$ws.Cells.Item(84, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(84, 16).Value = 1040
# Row 85
$ws.Cells.Item(85, 4).Value = 44167
$ws.Cells.Item(85, 10).Value = 110
# Row 86
$ws.Cells.Item(86, 4).Value = 44574
$ws.Cells.Item(86, 10).Value = 235
$ws.Cells.Item(86, 12).Value = 27000
$ws.Cells.Item(86, 13).Value = 26064
$ws.Cells.Item(86, 16).Value = 1043
# Row 87
$ws.Cells.Item(87, 4).Value = 44194
$ws.Cells.Item(87, 10).Value = 110
$ws.Cells.Item(87, 11).Value = 25000
$ws.Cells.Item(87, 12).Value = 25000
$ws.Cells.Item(87, 13).Value = 25000
$ws.Cells.Item(87, 16).Value = 1000
# Row 88
$ws.Cells.Item(88, 4).Value = 44516
$ws.Cells.Item(88, 10).Value = 35
$ws.Cells.Item(88, 15).Value = "Región Metropolitana"
# Row 89
$ws.Cells.Item(89, 4).Value = 44516
$ws.Cells.Item(89, 10).Value = 55
$ws.Cells.Item(89, 11).Value = 16000
$ws.Cells.Item(89, 12).Value = 16000
$ws.Cells.Item(89, 13).Value = 16000
$ws.Cells.Item(89, 15).Value = "Región del Maule"
$ws.Cells.Item(89, 16).Value = 640
# Row 90
$ws.Cells.Item(90, 4).Value = 44257
$ws.Cells.Item(90, 10).Value = 40
$ws.Cells.Item(90, 11).Value = 20000
$ws.Cells.Item(90, 12).Value = 20000
$ws.Cells.Item(90, 13).Value = 20000
$ws.Cells.Item(90, 16).Value = 800
# Row 91
$ws.Cells.Item(91, 4).Value = 44217
$ws.Cells.Item(91, 10).Value = 100
$ws.Cells.Item(91, 11).Value = 23000
$ws.Cells.Item(91, 12).Value = 23000
$ws.Cells.Item(91, 13).Value = 23000
$ws.Cells.Item(91, 16).Value = 920
# Row 92
$ws.Cells.Item(92, 4).Value = 44553
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 10).Value = 155
$ws.Cells.Item(92, 12).Value = 15000
$ws.Cells.Item(92, 13).Value = 15000
$ws.Cells.Item(92, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(92, 16).Value = 600
# Row 93
$ws.Cells.Item(93, 4).Value = 44523
$ws.Cells.Item(93, 10).Value = 20
$ws.Cells.Item(93, 11).Value = 15000
$ws.Cells.Item(93, 12).Value = 15000
$ws.Cells.Item(93, 13).Value = 15000
$ws.Cells.Item(93, 16).Value = 600
# Row 94
$ws.Cells.Item(94, 4).Value = 44519
$ws.Cells.Item(94, 8).Value = "Perfection"
$ws.Cells.Item(94, 10).Value = 65
$ws.Cells.Item(94, 11).Value = 15000
$ws.Cells.Item(94, 12).Value = 16000
$ws.Cells.Item(94, 13).Value = 15538
$ws.Cells.Item(94, 15).Value = "Región del Maule"
$ws.Cells.Item(94, 16).Value = 622
# Row 95
$ws.Cells.Item(95, 4).Value = 44242
$ws.Cells.Item(95, 10).Value = 85
$ws.Cells.Item(95, 11).Value = 17000
$ws.Cells.Item(95, 12).Value = 17000
$ws.Cells.Item(95, 13).Value = 17000
$ws.Cells.Item(95, 16).Value = 680
# Row 96
$ws.Cells.Item(96, 4).Value = 44568
$ws.Cells.Item(96, 10).Value = 50
$ws.Cells.Item(96, 11).Value = 18000
$ws.Cells.Item(96, 12).Value = 18000
$ws.Cells.Item(96, 13).Value = 18000
$ws.Cells.Item(96, 16).Value = 720
# Row 97
$ws.Cells.Item(97, 4).Value = 44529
$ws.Cells.Item(97, 10).Value = 110
$ws.Cells.Item(97, 11).Value = 15000
$ws.Cells.Item(97, 13).Value = 15000
$ws.Cells.Item(97, 16).Value = 600
# Row 98
$ws.Cells.Item(98, 4).Value = 44529
$ws.Cells.Item(98, 10).Value = 210
$ws.Cells.Item(98, 11).Value = 15000
$ws.Cells.Item(98, 12).Value = 15000
$ws.Cells.Item(98, 13).Value = 15000
$ws.Cells.Item(98, 15).Value = "Región del Maule"
$ws.Cells.Item(98, 16).Value = 600
# Row 99
$ws.Cells.Item(99, 4).Value = 44504
$ws.Cells.Item(99, 10).Value = 205
$ws.Cells.Item(99, 11).Value = 17000
$ws.Cells.Item(99, 13).Value = 17463
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 699
# Row 100
$ws.Cells.Item(100, 4).Value = 44498
$ws.Cells.Item(100, 10).Value = 20
$ws.Cells.Item(100, 11).Value = 20000
$ws.Cells.Item(100, 12).Value = 20000
$ws.Cells.Item(100, 13).Value = 20000
$ws.Cells.Item(100, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(100, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(100, 16).Value = 800
# Row 101
$ws.Cells.Item(101, 4).Value = 44616
$ws.Cells.Item(101, 10).Value = 50
$ws.Cells.Item(101, 11).Value = 23000
$ws.Cells.Item(101, 12).Value = 23000
$ws.Cells.Item(101, 13).Value = 23000
$ws.Cells.Item(101, 16).Value = 920
# Row 102
$ws.Cells.Item(102, 4).Value = 44560
$ws.Cells.Item(102, 10).Value = 65
$ws.Cells.Item(102, 11).Value = 15000
$ws.Cells.Item(102, 12).Value = 15000
$ws.Cells.Item(102, 13).Value = 15000
$ws.Cells.Item(102, 16).Value = 600
# Row 103
$ws.Cells.Item(103, 4).Value = 44901
$ws.Cells.Item(103, 10).Value = 30
$ws.Cells.Item(103, 11).Value = 25000
$ws.Cells.Item(103, 12).Value = 25000
$ws.Cells.Item(103, 13).Value = 25000
$ws.Cells.Item(103, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(103, 15).Value = "Región del Maule"
$ws.Cells.Item(103, 16).Value = 1000
# Row 104
$ws.Cells.Item(104, 4).Value = 44537
$ws.Cells.Item(104, 10).Value = 95
$ws.Cells.Item(104, 11).Value = 13000
$ws.Cells.Item(104, 12).Value = 13000
$ws.Cells.Item(104, 13).Value = 13000
$ws.Cells.Item(104, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(104, 16).Value = 520
# Row 105
$ws.Cells.Item(105, 4).Value = 44211
$ws.Cells.Item(105, 10).Value = 110
$ws.Cells.Item(105, 11).Value = 17000
$ws.Cells.Item(105, 12).Value = 17000
$ws.Cells.Item(105, 13).Value = 17000
$ws.Cells.Item(105, 16).Value = 680
# Row 106
$ws.Cells.Item(106, 4).Value = 44922
$ws.Cells.Item(106, 10).Value = 35
$ws.Cells.Item(106, 11).Value = 25000
$ws.Cells.Item(106, 12).Value = 25000
$ws.Cells.Item(106, 13).Value = 25000
$ws.Cells.Item(106, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(106, 16).Value = 1000
# Row 107
$ws.Cells.Item(107, 4).Value = 44897
$ws.Cells.Item(107, 10).Value = 20
$ws.Cells.Item(107, 11).Value = 26000
$ws.Cells.Item(107, 12).Value = 26000
$ws.Cells.Item(107, 13).Value = 26000
$ws.Cells.Item(107, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(107, 16).Value = 1040
# Row 108
$ws.Cells.Item(108, 4).Value = 44259
$ws.Cells.Item(108, 10).Value = 30
$ws.Cells.Item(108, 11).Value = 20000
$ws.Cells.Item(108, 12).Value = 20000
$ws.Cells.Item(108, 13).Value = 20000
$ws.Cells.Item(108, 16).Value = 800
# Row 109
$ws.Cells.Item(109, 4).Value = 44855
$ws.Cells.Item(109, 10).Value = 30
$ws.Cells.Item(109, 11).Value = 19000
$ws.Cells.Item(109, 13).Value = 19333
$ws.Cells.Item(109, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(109, 16).Value = 773
# Row 110
$ws.Cells.Item(110, 4).Value = 44522
$ws.Cells.Item(110, 10).Value = 80
$ws.Cells.Item(110, 11).Value = 16000
$ws.Cells.Item(110, 12).Value = 16000
$ws.Cells.Item(110, 13).Value = 16000
$ws.Cells.Item(110, 16).Value = 640
# Row 111
$ws.Cells.Item(111, 4).Value = 44196
$ws.Cells.Item(111, 10).Value = 80
$ws.Cells.Item(111, 11).Value = 14000
$ws.Cells.Item(111, 12).Value = 14000
$ws.Cells.Item(111, 13).Value = 14000
$ws.Cells.Item(111, 16).Value = 560
# Row 112
$ws.Cells.Item(112, 4).Value = 44161
$ws.Cells.Item(112, 10).Value = 100
$ws.Cells.Item(112, 11).Value = 23000
$ws.Cells.Item(112, 12).Value = 23000
$ws.Cells.Item(112, 13).Value = 23000
$ws.Cells.Item(112, 16).Value = 920
# Row 113
$ws.Cells.Item(113, 4).Value = 44161
$ws.Cells.Item(113, 10).Value = 100
$ws.Cells.Item(113, 11).Value = 20000
$ws.Cells.Item(113, 12).Value = 21000
$ws.Cells.Item(113, 13).Value = 20500
$ws.Cells.Item(113, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(113, 16).Value = 820
# Row 114
$ws.Cells.Item(114, 4).Value = 44201
$ws.Cells.Item(114, 10).Value = 30
$ws.Cells.Item(114, 11).Value = 18000
$ws.Cells.Item(114, 12).Value = 18000
$ws.Cells.Item(114, 13).Value = 18000
$ws.Cells.Item(114, 16).Value = 720
# Row 115
$ws.Cells.Item(115, 4).Value = 44600
$ws.Cells.Item(115, 10).Value = 10
$ws.Cells.Item(115, 11).Value = 22000
$ws.Cells.Item(115, 12).Value = 22000
$ws.Cells.Item(115, 13).Value = 22000
$ws.Cells.Item(115, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(115, 16).Value = 880
# Row 116
$ws.Cells.Item(116, 4).Value = 44484
$ws.Cells.Item(116, 10).Value = 10
$ws.Cells.Item(116, 11).Value = 25000
$ws.Cells.Item(116, 12).Value = 25000
$ws.Cells.Item(116, 13).Value = 25000
$ws.Cells.Item(116, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(116, 16).Value = 1000
# Row 117
$ws.Cells.Item(117, 4).Value = 45012
$ws.Cells.Item(117, 10).Value = 30
$ws.Cells.Item(117, 11).Value = 30000
$ws.Cells.Item(117, 12).Value = 30000
$ws.Cells.Item(117, 13).Value = 30000
$ws.Cells.Item(117, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(117, 16).Value = 1200
# Row 118
$ws.Cells.Item(118, 4).Value = 44222
$ws.Cells.Item(118, 10).Value = 65
$ws.Cells.Item(118, 11).Value = 18000
$ws.Cells.Item(118, 12).Value = 18000
$ws.Cells.Item(118, 13).Value = 18000
$ws.Cells.Item(118, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(118, 16).Value = 720
# Row 119
$ws.Cells.Item(119, 4).Value = 44540
$ws.Cells.Item(119, 10).Value = 120
$ws.Cells.Item(119, 11).Value = 15000
$ws.Cells.Item(119, 12).Value = 15000
$ws.Cells.Item(119, 13).Value = 15000
$ws.Cells.Item(119, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(119, 16).Value = 600
# Row 120
$ws.Cells.Item(120, 4).Value = 44567
$ws.Cells.Item(120, 10).Value = 30
$ws.Cells.Item(120, 11).Value = 18000
$ws.Cells.Item(120, 12).Value = 18000
$ws.Cells.Item(120, 13).Value = 18000
$ws.Cells.Item(120, 16).Value = 720
# Row 121
$ws.Cells.Item(121, 4).Value = 44525
$ws.Cells.Item(121, 10).Value = 100
$ws.Cells.Item(121, 11).Value = 16000
$ws.Cells.Item(121, 12).Value = 16000
$ws.Cells.Item(121, 13).Value = 16000
$ws.Cells.Item(121, 16).Value = 640
# Row 122
$ws.Cells.Item(122, 4).Value = 44557
$ws.Cells.Item(122, 10).Value = 375
$ws.Cells.Item(122, 11).Value = 10000
$ws.Cells.Item(122, 12).Value = 12000
$ws.Cells.Item(122, 13).Value = 10667
$ws.Cells.Item(122, 16).Value = 427
# Row 123
$ws.Cells.Item(123, 4).Value = 44200
# Row 124
$ws.Cells.Item(124, 4).Value = 44238
$ws.Cells.Item(124, 10).Value = 110
$ws.Cells.Item(124, 11).Value = 17000
$ws.Cells.Item(124, 12).Value = 17000
$ws.Cells.Item(124, 13).Value = 17000
$ws.Cells.Item(124, 16).Value = 680

Write-Output "Applied changes"